$wb = $excel.ActiveWorkbook

# Update the raw metric values on the "Metrics" sheet (B2:B13).
# All downstream formulas (today!B11:B22, E11:E22, F11:F22) reference these
# cells and will recalculate automatically.
$metrics = $wb.Worksheets.Item("Metrics")
$metrics.Range("B2").Value = 272268.55
$metrics.Range("B3").Value = 220449.92000000001
$metrics.Range("B4").Value = 86159.77
$metrics.Range("B5").Value = 10606
$metrics.Range("B6").Value = 4191519.4299999997
$metrics.Range("B7").Value = 3547977.3999999994
$metrics.Range("B8").Value = 1215525.4500000002
$metrics.Range("B9").Value = 161766
$metrics.Range("B10").Value = 32656843.230999827
$metrics.Range("B11").Value = 19577847.470000003
$metrics.Range("B12").Value = 11497234.34
$metrics.Range("B13").Value = 1259393

# Move the selection on the Metrics sheet.
$metrics.Range("E9").Select() | Out-Null

# Move the selection on the "today" sheet (the active sheet in the workbook).
$today = $wb.Worksheets.Item("today")
$today.Range("F7").Select() | Out-Null
